# A new weekly price record (Albahaca, Feria Lagunitas de Puerto Montt) needs
# to be inserted as the new most-recent entry, right above the existing row 43.
# Inserting an entire row shifts all the following rows (old 43..171) down by
# one (to 44..172), which matches the target diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("43:43").Insert()

$ws.Cells.Item(43, 1).Value2  = 4
$ws.Cells.Item(43, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(43, 3).Value2  = "Los Lagos"
$ws.Cells.Item(43, 4).Value2  = 44980
$ws.Cells.Item(43, 5).Value2  = 10
$ws.Cells.Item(43, 6).Value2  = 100112052
$ws.Cells.Item(43, 7).Value2  = "Albahaca"
$ws.Cells.Item(43, 8).Value2  = "Sin especificar"
$ws.Cells.Item(43, 9).Value2  = "Primera"
$ws.Cells.Item(43, 10).Value2 = 40
$ws.Cells.Item(43, 11).Value2 = 5500
$ws.Cells.Item(43, 12).Value2 = 6000
$ws.Cells.Item(43, 13).Value2 = 5750
$ws.Cells.Item(43, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(43, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(43, 16).Value2 = 958
$ws.Cells.Item(43, 17).Value2 = 6
$ws.Cells.Item(43, 18).Value2 = "Hortaliza"
